# Estadisticos Matutinos 15 Oct
$wb = $excel.ActiveWorkbook

# --- Sheet: Estadisticos 1P ---
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")
$ws1.Range("D6").Value = 4
$ws1.Range("F6").Value = 20
$ws1.Range("G6").Value = 83.33
$ws1.Range("H6").Value = 8.800000000000001

# --- Sheet: Estadisticos 2P ---
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Range("E6").Value = 20

# --- Sheet: Estadisticos Final ---
$ws3 = $wb.Worksheets.Item("Estadisticos Final")
$ws3.Range("D6").Value = 4
$ws3.Range("F6").Value = 20
$ws3.Range("G6").Value = 83.33
$ws3.Range("H6").Value = 8.800000000000001

# --- Sheet: Rescatables ---
$ws4 = $wb.Worksheets.Item("Rescatables")

$colA = @(21330051920199, 21330051920201, 21330051920306, 21330051920328, 21330051920330, 21330051920352, 21330051920283, 21330051920311)
$colB = @("MORALES", "OFICIAL", "LOPEZ", "TRUJILLO", "ZUÑIGA", "SANCHEZ", "TELLEZ", "OFICIAL")
$colC = @("ARELLANO", "TZOMPAXTLE", "RAMOS", "OSORIO", "ESPINDOLA", "CASTELLANOS", "VALENCIA", "VILLASEÑOR")
$colD = @("WENCESLAO", "CLEMENTE", "ANETTE JOCELYN", "KARINA YOSELIN", "ROSA ITZEL", "NOEMI", "VIANEY", "MONICA AIME")
$colE = @("QUÍMICA I", "QUÍMICA I", "QUÍMICA I", "QUÍMICA I", "QUÍMICA I", "QUÍMICA I", "QUÍMICA I", "QUÍMICA I")
$colF = @("1BM", "1BM", "1EM", "1EM", "1EM", "1FM", "1DM", "1EM")
$colG = @(6, 6, 6, 6, 6, 6, 6, 6)

for ($i = 0; $i -lt $colA.Length; $i++) {
    $ws4.Cells.Item($i + 2, 1).Value = $colA[$i]
}
for ($i = 0; $i -lt $colB.Length; $i++) {
    $ws4.Cells.Item($i + 2, 2).Value = $colB[$i]
}
for ($i = 0; $i -lt $colC.Length; $i++) {
    $ws4.Cells.Item($i + 2, 3).Value = $colC[$i]
}
for ($i = 0; $i -lt $colD.Length; $i++) {
    $ws4.Cells.Item($i + 2, 4).Value = $colD[$i]
}
for ($i = 0; $i -lt $colE.Length; $i++) {
    $ws4.Cells.Item($i + 2, 5).Value = $colE[$i]
}
for ($i = 0; $i -lt $colF.Length; $i++) {
    $ws4.Cells.Item($i + 2, 6).Value = $colF[$i]
}
for ($i = 0; $i -lt $colG.Length; $i++) {
    $ws4.Cells.Item($i + 2, 7).Value = $colG[$i]
}
